# Sprint 4 User Stories workbook update:
#  - B1 header text fixed to "User stories" (was a mangled "User stories<TAB>Acceptance Criteria")
#  - New "Estimation" column (D) added with story-point values for every existing row
#  - Two new user stories appended (FNDJL-22, FNDJL-23) with story / acceptance criteria / estimate
#  - B6:B7 acceptance-criteria cells normalized to the same wrap/top style as the rest of column B/C
#  - Selection / top-left-cell view state updated to reflect the new bottom-of-sheet position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): fix B1 text, add D1 "Estimation" header
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "User stories"

$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D1").Value = "Estimation"

# Header cells all get vertical-top alignment on top of their existing fill/border
$ws.Range("A1:D1").VerticalAlignment = -4160   # xlVAlignTop

# ---------------------------------------------------------------------------
# 2. Normalize B6:B7 to the same wrap+top style used across the rest of B/C
# ---------------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("B6:B7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. New column D (Estimation) for existing rows 2-7
#    (work out the centered/no-wrap format on D2 once, then fan it out so we
#    don't leave a trail of transient single-axis-alignment styles behind)
# ---------------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").WrapText = $false
$ws.Range("D2").HorizontalAlignment = -4108      # xlHAlignCenter
$ws.Range("D2").VerticalAlignment = -4108        # xlVAlignCenter

$ws.Range("D2").Copy()
$ws.Range("D3:D7").PasteSpecial(-4122)

$ws.Range("D2").Value = 13
$ws.Range("D3").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 8
$ws.Range("D6").Value = 8
$ws.Range("D7").Value = 3

# ---------------------------------------------------------------------------
# 4. New rows 8 and 9 (FNDJL-22, FNDJL-23)
# ---------------------------------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("B8:C9").PasteSpecial(-4122)

$ws.Range("D7").Copy()
$ws.Range("D8:D9").PasteSpecial(-4122)

$ws.Rows.Item(8).RowHeight = 102
$ws.Rows.Item(9).RowHeight = 102

$ws.Range("A8").Value = "FNDJL-22"
$ws.Range("B8").Value = "As a business owner `nI should able to upload post related to news`nSo that I can make sure that news is correct through users feedback"
$ws.Range("C8").Value = "User should be able to upload text blog to a portal and other users should be able to read it."
$ws.Range("D8").Value = 5

$ws.Range("A9").Value = "FNDJL-23"
$ws.Range("B9").Value = "As a Business Owner `nI want to comment on blog post created by other users`nSo that I can validate the news posted on blog`n"
$ws.Range("C9").Value = "User must be able to write a text comment on blogs posted by other users and it should be visible to other users"
$ws.Range("D9").Value = 4

# ---------------------------------------------------------------------------
# 5. View state: scroll so row 5 is at top, select C8
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("C8").Select()
